# Adjustment sensitivity EOL RIR
# Update the "2020" scenario-year column header to "2030" on every sheet,
# and refresh the corresponding sensitivity values (column C) for the
# sheets where the underlying numbers changed (Neodymium, Copper, Raw silicon).
# Dysprosium's column C values are unchanged (all zero).

$wb = $excel.ActiveWorkbook

function Update-Sheet {
    param(
        [string]$SheetName,
        [string[]]$Values   # $null => only update the header
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Header: year 2020 -> 2030
    $ws.Range("C1").Value = 2030

    if ($Values -and $Values.Count -gt 0) {
        $ws.Range("C2").Value = [double]$Values[0]
        $ws.Range("C3").Value = [double]$Values[1]
        $ws.Range("C4").Value = [double]$Values[2]
        $ws.Range("C5").Value = [double]$Values[3]
    }
}

Update-Sheet "Neodymium"  @("2.195405251500087E-06", "0.0001062411525673284", "9.608716352691787E-05", "2.138791829054013E-09")
Update-Sheet "Dysprosium" $null
Update-Sheet "Copper"     @("0.003816340722347758", "0.0137679456486295", "0.003685389348936462", "0.008080293662215502")
Update-Sheet "Raw silicon" @("5.750015024097242E-05", "0.0001921210602835477", "5.3924808017845E-05", "6.847896595910317E-05")
